$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.4002345243970353
$ws.Range("C2").Value = 0.02923790351773903
$ws.Range("D2").Value = 0.1684952640468538
$ws.Range("E2").Value = 0.1558357427434167
$ws.Range("F2").Value = 1.589051218795767
$ws.Range("J2").Value = 0.1840581092675535
$ws.Range("K2").Value = 0.3524572108632924
$ws.Range("M2").Value = 0.2310193804042413
$ws.Range("N2").Value = 1.941613291518943
$ws.Range("O2").Value = 3.949852420747902

$ws.Range("B3").Value = 0.3678260096117185
$ws.Range("C3").Value = 0.02554740063324346
$ws.Range("D3").Value = 0.1653073817915782
$ws.Range("E3").Value = 0.154639316785449
$ws.Range("F3").Value = 1.592261207634138
$ws.Range("J3").Value = 0.1838083673130768
$ws.Range("K3").Value = 0.3189940532912203
$ws.Range("M3").Value = 0.2208632601735161
$ws.Range("N3").Value = 1.960455524992244
$ws.Range("O3").Value = 3.966629397327239

$ws.Range("B4").Value = 0.348042742285628
$ws.Range("C4").Value = 0.02327198865255298
$ws.Range("D4").Value = 0.1634166907768702
$ws.Range("E4").Value = 0.1539715075598131
$ws.Range("F4").Value = 1.59496774858335
$ws.Range("J4").Value = 0.1837366976299393
$ws.Range("K4").Value = 0.2985130855273468
$ws.Range("M4").Value = 0.2147240888940054
$ws.Range("N4").Value = 1.972615078615124
$ws.Range("O4").Value = 3.978778437857088

$ws.Range("B5").Value = 0.3400104637549646
$ws.Range("C5").Value = 0.02234241728484676
$ws.Range("D5").Value = 0.1626630749417899
$ws.Range("E5").Value = 0.1537162120047739
$ws.Range("F5").Value = 1.596255840223357
$ws.Range("J5").Value = 0.1837280565282029
$ws.Range("K5").Value = 0.2901838747440024
$ws.Range("M5").Value = 0.2122468196736911
$ws.Range("N5").Value = 1.977718665748937
$ws.Range("O5").Value = 3.98419421570226

$ws.Range("B6").Value = 0.338678509479081
$ws.Range("C6").Value = 0.02218792381748358
$ws.Range("D6").Value = 0.1625389582756469
$ws.Range("E6").Value = 0.1536748389023401
$ws.Range("F6").Value = 1.596480915109083
$ws.Range("J6").Value = 0.1837278645595859
$ws.Range("K6").Value = 0.2888018531401286
$ws.Range("M6").Value = 0.2118369553211537
$ws.Range("N6").Value = 1.978575080459959
$ws.Range("O6").Value = 3.985121591112929

$ws.Range("B7").Value = 0.3479342958412701
$ws.Range("C7").Value = 0.02325946146797264
$ws.Range("D7").Value = 0.1634064588824202
$ws.Range("E7").Value = 0.1539679963023488
$ws.Range("F7").Value = 1.594984370313128
$ws.Range("J7").Value = 0.183736497789539
$ws.Range("K7").Value = 0.2984006856286925
$ws.Range("M7").Value = 0.2146905801785195
$ws.Range("N7").Value = 1.9726833061952
$ws.Range("O7").Value = 3.978849594025604

$ws.Range("B8").Value = 0.3890363294032682
$ws.Range("C8").Value = 0.02796740286041199
$ws.Range("D8").Value = 0.1673822942558871
$ws.Range("E8").Value = 0.1554093764510824
$ws.Range("F8").Value = 1.590005457046161
$ws.Range("J8").Value = 0.183955063070151
$ws.Range("K8").Value = 0.3409057801355004
$ws.Range("M8").Value = 0.2274975727689679
$ws.Range("N8").Value = 1.94798758287069
$ws.Range("O8").Value = 3.955253747855977

$ws.Range("B9").Value = 0.4705386873024793
$ws.Range("C9").Value = 0.03712315462128402
$ws.Range("D9").Value = 0.1757046279796555
$ws.Range("E9").Value = 0.1587643007657995
$ws.Range("F9").Value = 1.586071765778328
$ws.Range("J9").Value = 0.185030718789605
$ws.Range("K9").Value = 0.4247618543170972
$ws.Range("M9").Value = 0.2533735112605697
$ws.Range("N9").Value = 1.904241703089235
$ws.Range("O9").Value = 3.92363511658138

$ws.Range("B10").Value = 0.5309518586556692
$ws.Range("C10").Value = 0.04380158918951338
$ws.Range("D10").Value = 0.1821358201214025
$ws.Range("E10").Value = 0.1615494938844044
$ws.Range("F10").Value = 1.586728821723739
$ws.Range("J10").Value = 0.1862144641841397
$ws.Range("K10").Value = 0.4866621576779266
$ws.Range("M10").Value = 0.2728428026373919
$ws.Range("N10").Value = 1.874952586932143
$ws.Range("O10").Value = 3.90932835186797

$ws.Range("B11").Value = 0.5585481322935095
$ws.Range("C11").Value = 0.04682898100931254
$ws.Range("D11").Value = 0.1851295761479719
$ws.Range("E11").Value = 0.1628857550421579
$ws.Range("F11").Value = 1.587796632561066
$ws.Range("J11").Value = 0.1868382272289253
$ws.Range("K11").Value = 0.5148824520861126
$ws.Range("M11").Value = 0.2817982479464192
$ws.Range("N11").Value = 1.862246525303703
$ws.Range("O11").Value = 3.904755804067861

$ws.Range("B12").Value = 0.5690141187568827
$ws.Range("C12").Value = 0.04797380137172524
$ws.Range("D12").Value = 0.1862729574668691
$ws.Range("E12").Value = 0.1634016835525749
$ws.Range("F12").Value = 1.588311408281854
$ws.Range("J12").Value = 0.187086669728437
$ws.Range("K12").Value = 0.5255771964204143
$ws.Range("M12").Value = 0.2852035000611934
$ws.Range("N12").Value = 1.857523878761639
$ws.Range("O12").Value = 3.903302457282109

$ws.Range("B13").Value = 0.5667593834730269
$ws.Range("C13").Value = 0.04772731527637575
$ws.Range("D13").Value = 0.1860262794221796
$ws.Range("E13").Value = 0.1632901288358077
$ws.Range("F13").Value = 1.588195633595532
$ws.Range("J13").Value = 0.1870326193460699
$ws.Range("K13").Value = 0.5232735295835766
$ws.Range("M13").Value = 0.2844694974982929
$ws.Range("N13").Value = 1.858537032341656
$ws.Range("O13").Value = 3.903603092176468

$ws.Range("B14").Value = 0.5594088600862221
$ws.Range("C14").Value = 0.04692319818987301
$ws.Range("D14").Value = 0.1852234487266742
$ws.Range("E14").Value = 0.1629280023268684
$ws.Range("F14").Value = 1.587836771368302
$ws.Range("J14").Value = 0.1868584216385401
$ws.Range("K14").Value = 0.5157621511336004
$ws.Range("M14").Value = 0.2820781201698352
$ws.Range("N14").Value = 1.861856209151052
$ws.Range("O14").Value = 3.904630662224463

$ws.Range("B15").Value = 0.554908504442011
$ws.Range("C15").Value = 0.04643044494785897
$ws.Range("D15").Value = 0.1847329536414719
$ws.Range("E15").Value = 0.1627074793407886
$ws.Range("F15").Value = 1.587631333247458
$ws.Range("J15").Value = 0.1867533133869941
$ws.Range("K15").Value = 0.5111622845591626
$ws.Range("M15").Value = 0.2806151527934233
$ws.Range("N15").Value = 1.863900874579434
$ws.Range("O15").Value = 3.905296300491443

$ws.Range("B16").Value = 0.5291506272581614
$ws.Range("C16").Value = 0.04360352285118552
$ws.Range("D16").Value = 0.1819415350629896
$ws.Range("E16").Value = 0.1614635560184574
$ws.Range("F16").Value = 1.586674497464969
$ws.Range("J16").Value = 0.1861754135485754
$ws.Range("K16").Value = 0.4848190891079298
$ws.Range("M16").Value = 0.2722595150139
$ws.Range("N16").Value = 1.875795390927875
$ws.Range("O16").Value = 3.909666114793794

$ws.Range("B17").Value = 0.5133778376961402
$ws.Range("C17").Value = 0.04186652852810369
$ws.Range("D17").Value = 0.1802464854931145
$ws.Range("E17").Value = 0.160718157589951
$ws.Range("F17").Value = 1.586284340195675
$ws.Range("J17").Value = 0.1858427130199445
$ws.Range("K17").Value = 0.4686737847660822
$ws.Range("M17").Value = 0.2671587682059311
$ws.Range("N17").Value = 1.883250526598497
$ws.Range("O17").Value = 3.91284250561381

$ws.Range("B18").Value = 0.5043165094060953
$ws.Range("C18").Value = 0.04086645505169884
$ws.Range("D18").Value = 0.1792779602543675
$ws.Range("E18").Value = 0.160295946141904
$ws.Range("F18").Value = 1.586132314727166
$ws.Range("J18").Value = 0.1856593801114315
$ws.Range("K18").Value = 0.4593932557597498
$ws.Range("M18").Value = 0.264234260946445
$ws.Range("N18").Value = 1.887596662243878
$ws.Range("O18").Value = 3.914851686231628

$ws.Range("B19").Value = 0.501250363722221
$ws.Range("C19").Value = 0.04052767720364159
$ws.Range("D19").Value = 0.1789511401285893
$ws.Range("E19").Value = 0.1601541141165583
$ws.Range("F19").Value = 1.586093277873005
$ws.Range("J19").Value = 0.1855986862522059
$ws.Range("K19").Value = 0.4562520439684192
$ws.Range("M19").Value = 0.2632456771700546
$ws.Range("N19").Value = 1.889078173374518
$ws.Range("O19").Value = 3.915563259257823

$ws.Range("B20").Value = 0.5150557678906864
$ws.Range("C20").Value = 0.04205153845629184
$ws.Range("D20").Value = 0.1804262623854953
$ws.Range("E20").Value = 0.1607968318436406
$ws.Range("F20").Value = 1.586318382851516
$ws.Range("J20").Value = 0.185877298909638
$ws.Range("K20").Value = 0.4703918810614596
$ws.Range("M20").Value = 0.2677007894574288
$ws.Range("N20").Value = 1.882450897263004
$ws.Range("O20").Value = 3.912485517116608

$ws.Range("B21").Value = 0.5615674602411787
$ws.Range("C21").Value = 0.04715943037781756
$ws.Range("D21").Value = 0.1854589967501568
$ws.Range("E21").Value = 0.1630340989143129
$ws.Range("F21").Value = 1.587939182360145
$ws.Range("J21").Value = 0.1869092558190886
$ws.Range("K21").Value = 0.5179682030806703
$ws.Range("M21").Value = 0.2827801466404196
$ws.Range("N21").Value = 1.860878874424428
$ws.Range("O21").Value = 3.904321291674819

$ws.Range("B22").Value = 0.5920578391317122
$ws.Range("C22").Value = 0.05048845642076571
$ws.Range("D22").Value = 0.1888047390519318
$ws.Range("E22").Value = 0.1645540561624692
$ws.Range("F22").Value = 1.58964200491755
$ws.Range("J22").Value = 0.1876550141396507
$ws.Range("K22").Value = 0.549110400445727
$ws.Range("M22").Value = 0.2927170083948667
$ws.Range("N22").Value = 1.847298292990331
$ws.Range("O22").Value = 3.900606853225469

$ws.Range("B23").Value = 0.5757762764706342
$ws.Range("C23").Value = 0.04871255924835793
$ws.Range("D23").Value = 0.1870139088177751
$ws.Range("E23").Value = 0.1637375547856728
$ws.Range("F23").Value = 1.588674342392991
$ws.Range("J23").Value = 0.1872504718325132
$ws.Range("K23").Value = 0.5324849734141992
$ws.Range("M23").Value = 0.2874061097376241
$ws.Range("N23").Value = 1.854499096173436
$ws.Range("O23").Value = 3.902441020758346

$ws.Range("B24").Value = 0.5142971548700643
$ws.Range("C24").Value = 0.04196789998734118
$ws.Range("D24").Value = 0.1803449666153227
$ws.Range("E24").Value = 0.1607612434994579
$ws.Range("F24").Value = 1.586302767004312
$ws.Range("J24").Value = 0.1858616378918541
$ws.Range("K24").Value = 0.4696151246368743
$ws.Range("M24").Value = 0.2674557167794305
$ws.Range("N24").Value = 1.882812222543554
$ws.Range("O24").Value = 3.912646341469554

$ws.Range("B25").Value = 0.4483953042795861
$ws.Range("C25").Value = 0.03465464304153443
$ws.Range("D25").Value = 0.1733973179347004
$ws.Range("E25").Value = 0.1578003222436521
$ws.Range("F25").Value = 1.586512516800305
$ws.Range("J25").Value = 0.1846705573212404
$ws.Range("K25").Value = 0.4020242734504222
$ws.Range("M25").Value = 0.2462924803521531
$ws.Range("N25").Value = 1.915575217445694
$ws.Range("O25").Value = 3.930621077509613
